{"js": "// Replace the payment-provider name \"braintree\" with \"stripe\" in the body\n// of the document (commit: \"braintree is now stripe\").\nconst body = context.document.body;\n\nconst results = body.search(\"braintree\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"stripe\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the payment-provider name \"braintree\" with \"stripe\" in the body\n# of the document (commit: \"braintree is now stripe\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"braintree\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"stripe\"\n\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
